# Insert a new weekly price-report row for "Apio" (Vega Modelo de Temuco)
# right above the current row 512. Excel's Rows.Item(512).Insert() shifts
# row 512 (and everything below it, down through 618) down by one,
# automatically growing the sheet's used range from A1:R618 to A1:R619.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(512).Insert()

$ws.Range("A512").Value = 10
$ws.Range("B512").Value = "Vega Modelo de Temuco"
$ws.Range("C512").Value = "La Araucanía"
$ws.Range("D512").Value = 45275
$ws.Range("E512").Value = 9
$ws.Range("F512").Value = 100112017
$ws.Range("G512").Value = "Apio"
$ws.Range("H512").Value = "Americana (o)"
$ws.Range("I512").Value = "Primera"
$ws.Range("J512").Value = 110
$ws.Range("K512").Value = 10000
$ws.Range("L512").Value = 10000
$ws.Range("M512").Value = 10000
$ws.Range("N512").Value = "$/caja 8 unidades"
$ws.Range("O512").Value = "Provincia del Elquí"
$ws.Range("P512").Value = 10000
$ws.Range("Q512").Value = 1
$ws.Range("R512").Value = "Hortaliza"
